$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the 11 newly-added medication names to the bottom of the list ---
$newDrugs = @(
    "바라크루드",
    "로수바미브",
    "트레시바",
    "안플라그",
    "하이드린",
    "도베실산칼슘수화물",
    "엔시움",
    "비리어드",
    "아토젯정",
    "아펜팔정",
    "아모잘탄정"
)

$startRow = 160
for ($i = 0; $i -lt $newDrugs.Length; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $newDrugs[$i]
    $cell.Style = $ws.Cells.Item($r - 1, 1).Style
}

# --- Re-sort the full data range (A2:A170), same as the existing autofilter sort state ---
$dataRange = $ws.Range("A2:A170")
$keyRange = $ws.Range("A2")
$dataRange.Sort($keyRange, 1, $null, $null, 2)

# --- Misc workbook-level metadata touched by this save ---
$wb.Sheets.Item(1).Select()
